# Apply updated crypto price/volume figures (and the TrustWalletToken / TheSandbox
# row 41/42 swap) as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.159.10"
$ws.Range("E2").Value = "'  +0.11%  "
$ws.Range("D3").Value = "'1.599.96"
$ws.Range("E3").Value = "'  -0.40%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "'  +0.22%  "
$ws.Range("E5").Value = "'  +0.17%  "
$ws.Range("D6").Value = "'302.62"
$ws.Range("E6").Value = "'  +0.12%  "
$ws.Range("D7").Value = "'0.3778"
$ws.Range("E7").Value = "'  -0.24%  "
$ws.Range("D8").Value = "'51.88"
$ws.Range("E8").Value = "'  +3.31%  "
$ws.Range("D9").Value = "'0.3614"
$ws.Range("E9").Value = "'  -1.51%  "
$ws.Range("D10").Value = "'1.262"
$ws.Range("E10").Value = "'  -1.08%  "
$ws.Range("E11").Value = "'  +0.13%  "
$ws.Range("D12").Value = "'0.08105"
$ws.Range("E12").Value = "'  -0.74%  "
$ws.Range("D13").Value = "'22.68"
$ws.Range("E13").Value = "'  -1.23%  "
$ws.Range("D14").Value = "'6.575"
$ws.Range("E14").Value = "'  -0.95%  "
$ws.Range("D15").Value = "'7.386"
$ws.Range("E15").Value = "'  -0.62%  "
$ws.Range("D16").Value = "'0.00001242"
$ws.Range("E16").Value = "'  -1.80%  "
$ws.Range("D17").Value = "'1.600.43"
$ws.Range("E17").Value = "'  -0.29%  "
$ws.Range("D18").Value = "'93.90"
$ws.Range("E18").Value = "'  +1.80%  "
$ws.Range("D19").Value = "'0.06889"
$ws.Range("E19").Value = "'  +0.15%  "
$ws.Range("D20").Value = "'18.03"
$ws.Range("E20").Value = "'  -1.69%  "
$ws.Range("D21").Value = "'6.529"
$ws.Range("E21").Value = "'  -1.26%  "
$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "'  +0.14%  "
$ws.Range("D23").Value = "'12.94"
$ws.Range("E23").Value = "'  -1.28%  "
$ws.Range("D24").Value = "'23.165.25"
$ws.Range("E24").Value = "'  +0.10%  "
$ws.Range("D25").Value = "'2.401"
$ws.Range("E25").Value = "'  +1.23%  "
$ws.Range("D26").Value = "'2.980"
$ws.Range("E26").Value = "'  +5.30%  "
$ws.Range("D27").Value = "'21.18"
$ws.Range("E27").Value = "'  -0.01%  "
$ws.Range("D28").Value = "'149.76"
$ws.Range("E28").Value = "'  -0.26%  "
$ws.Range("D29").Value = "'5.244"
$ws.Range("E29").Value = "'  -0.81%  "
$ws.Range("D30").Value = "'133.55"
$ws.Range("E30").Value = "'  -0.62%  "
$ws.Range("D31").Value = "'2.369"
$ws.Range("E31").Value = "'  -0.72%  "
$ws.Range("D32").Value = "'6.731"
$ws.Range("E32").Value = "'  -2.87%  "
$ws.Range("D33").Value = "'1.777.23"
$ws.Range("E33").Value = "'  -0.48%  "
$ws.Range("D34").Value = "'0.9650"
$ws.Range("E34").Value = "'  +0.11%  "
$ws.Range("D35").Value = "'0.07454"
$ws.Range("E35").Value = "'  -3.69%  "
$ws.Range("D36").Value = "'10.23"
$ws.Range("E36").Value = "'  -2.25%  "
$ws.Range("D37").Value = "'0.02710"
$ws.Range("E37").Value = "'  -0.89%  "
$ws.Range("D38").Value = "'0.2507"
$ws.Range("E38").Value = "'  -2.12%  "
$ws.Range("D39").Value = "'0.08803"
$ws.Range("D40").Value = "'6.043"
$ws.Range("E40").Value = "'  -4.43%  "
$ws.Range("B41").Value = "'TrustWalletToken"
$ws.Range("C41").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'1.359"
$ws.Range("E41").Value = "'  -0.86%  "
$ws.Range("B42").Value = "'TheSandbox"
$ws.Range("C42").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.7082"
$ws.Range("E42").Value = "'  -0.55%  "
$ws.Range("D43").Value = "'12.42"
$ws.Range("E43").Value = "'  -2.16%  "
$ws.Range("D44").Value = "'15.61"
$ws.Range("E44").Value = "'  +1.78%  "
$ws.Range("D45").Value = "'0.6517"
$ws.Range("E45").Value = "'  -2.04%  "
$ws.Range("D46").Value = "'2.303"
$ws.Range("E46").Value = "'  -0.97%  "
$ws.Range("D47").Value = "'4.012"
$ws.Range("E47").Value = "'  +0.08%  "
$ws.Range("D48").Value = "'131.98"
$ws.Range("E48").Value = "'  -0.45%  "
$ws.Range("D49").Value = "'0.07945"
$ws.Range("E49").Value = "'  -0.07%  "
$ws.Range("D50").Value = "'1.198"
$ws.Range("E50").Value = "'  -4.07%  "
$ws.Range("D51").Value = "'1.199"
$ws.Range("E51").Value = "'  -1.06%  "
